$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Text edits inside the big paragraph (originally Paragraph 3)
# ---------------------------------------------------------------------

# Remove the "One of the factors..." sentence from its original spot;
# it will be re-inserted later as its own paragraph near the end.
$r = $d.Content
$r.Find.Execute("One of the factors that economist study is related with the need of China to maintain surpluses in the balance trade. China requires a continuous growth in its exports, in order to secure growing labor stability. ", $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# Insert "the US or " before "China exposes smaller economies"
$r = $d.Content
$r.Find.Execute("Opening a free trade agreement with China exposes", $false, $false, $false, $false, $false, $true, 1, $false, "Opening a free trade agreement with the US or China exposes", 2) | Out-Null

# Replace "to balance the commercial power." with "to get a better-balanced commercial agreement."
$r = $d.Content
$r.Find.Execute("to balance the commercial power.", $false, $false, $false, $false, $false, $true, 1, $false, "to get a better-balanced commercial agreement.", 2) | Out-Null

# Insert "trend in its trade " before "balance (2)"
$r = $d.Content
$r.Find.Execute("improve the negative current commercial balance (2).", $false, $false, $false, $false, $false, $true, 1, $false, "improve the negative trend in its trade balance (2).", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Split the (former) big paragraph into two paragraphs, with a blank
#    paragraph in between, matching the new document structure.
# ---------------------------------------------------------------------

# Split right after "...get a better-balanced commercial agreement. "
# -> end of paragraph A here, followed by a blank paragraph.
$r = $d.Content
$r.Find.Execute("to get a better-balanced commercial agreement. ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertParagraphAfter()

Write-Output "--- after split 1 ---"
for ($i = 3; $i -le 6; $i++) {
    Write-Output "$i : [$($d.Paragraphs.Item($i).Range.Text)]"
}

# ---------------------------------------------------------------------
# 3) Insert a new blank paragraph plus a new paragraph (moved sentence)
#    right before the "References:" paragraph.
# ---------------------------------------------------------------------

function Get-ParaIndex($doc, $matchText) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text -like "*$matchText*") {
            return $i
        }
    }
    return -1
}

$refIdx = Get-ParaIndex $d "References:"
$refPara = $d.Paragraphs.Item($refIdx)
$r = $refPara.Range
$r.Collapse(1)
$r.InsertParagraphBefore()
$r.InsertParagraphBefore()

$refIdx2 = Get-ParaIndex $d "References:"
$newTextPara = $d.Paragraphs.Item($refIdx2 - 1)
$newTextPara.Range.InsertBefore("One of the factors that economist study is related with the need of China to maintain surpluses in the balance trade. China requires a continuous growth in its exports, in order to secure growing labor stability.")

Write-Output "--- paragraphs around insert ---"
for ($i = 4; $i -le 9; $i++) {
    Write-Output "$i : [$($d.Paragraphs.Item($i).Range.Text)]"
}

# ---------------------------------------------------------------------
# 4) Move the "_GoBack" bookmark from the old big paragraph to the end
#    of the paragraph holding the Narins_2018 hyperlink reference.
# ---------------------------------------------------------------------

$narinsIdx = Get-ParaIndex $d "Narins_2018_PG_Chinese_Trade.pdf"
$narinsPara = $d.Paragraphs.Item($narinsIdx)
$bmRange = $narinsPara.Range
$bmRange.MoveEnd(1, -1) | Out-Null
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
